$d = $word.ActiveDocument

# 1) Update the letter date: September 19, 2025 -> September 21, 2025
$range = $d.Content
$range.Find.Execute("September 19, 2025", $true, $false, $false, $false, $false, $true, 1, $false, "September 21, 2025", 2)

# 2) Split "7 Ludina Way, Redwood City CA 94061" (the mailing-address paragraph,
#    not the copy inside the property-address table) into two paragraphs:
#    "7 Ludina Way" and "Redwood City, CA 94061"
foreach ($p in $d.Paragraphs) {
    $r = $p.Range
    if ($r.Text.TrimEnd([char]13, [char]7) -eq "7 Ludina Way, Redwood City CA 94061") {
        $inTable = $false
        foreach ($t in $d.Tables) {
            if ($r.Start -ge $t.Range.Start -and $r.End -le $t.Range.End) {
                $inTable = $true
            }
        }
        if (-not $inTable) {
            # Replace the whole paragraph's text with the first line
            # ("7 Ludina Way"), keeping its trailing paragraph mark, then
            # insert a fresh paragraph mark right after it and fill the
            # newly-created (now empty) paragraph with the second line
            # ("Redwood City, CA 94061"). The paragraph that originally
            # followed (already empty in the source) is left untouched.
            $r.Text = "7 Ludina Way"
            $r.InsertParagraphAfter()
            $addrPara = $p.Next()
            $insertPoint = $d.Range($addrPara.Range.Start, $addrPara.Range.Start)
            $insertPoint.InsertAfter("Redwood City, CA 94061")
            break
        }
    }
}

# 3) Remove the empty "NoSpacing" paragraph right after
#    "Kentfield Pacific Place Board of Directors"
foreach ($p in $d.Paragraphs) {
    $r = $p.Range
    $txt = $r.Text.TrimEnd([char]13, [char]7)
    if ($txt -eq "Kentfield Pacific Place Board of Directors" -or $txt -eq "Board of Directors") {
        $next = $p.Next()
        if ($next -ne $null) {
            $nextTxt = $next.Range.Text.TrimEnd([char]13, [char]7)
            if ($nextTxt -eq "" -and $next.Style.NameLocal -eq "No Spacing") {
                $next.Range.Delete()
            }
        }
        break
    }
}
